$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

# Write the date as literal text (matching the sheet's existing convention of
# storing dates as plain strings, not real date serials). Using .Value
# directly on the target cell lets Excel's autodetection turn "08/13/2025"
# into a date serial + a brand new date-formatted style, which we don't want.
# Instead, stage it on a scratch cell pre-formatted as Text, then copy only
# the resulting value/type (not the format) onto A40, so A40 keeps the
# workbook's default (unstyled) cell formatting just like the rows above it.
$scratch = $ws.Cells.Item(1048576, 16384)
$scratch.NumberFormat = "@"
$scratch.Value = "08/13/2025"
$scratch.Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 2).Value = 0.0004084099999999979
$ws.Cells.Item($row, 3).Value = 122425.9934869378
$ws.Cells.Item($row, 4).Value = 50
